# Fiumalbo.xlsx - "aggiornamento fino a 27/05"
# Append 14 new daily rows (256-269) covering 2021-05-14 .. 2021-05-27
# (Excel date serials 44330-44343), each with B/C/D = 0, continuing the
# existing layout: column A keeps the bordered/bold/date-time style used
# by the rest of the date column, columns B/C/D stay unstyled.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 256
$startSerial = 44330
$numRows = 14

for ($i = 0; $i -lt $numRows; $i++) {
    $row = $startRow + $i
    $serial = $startSerial + $i

    $prevRow = $row - 1
    $src = $ws.Range("A" + $prevRow + ":D" + $prevRow)
    $dst = $ws.Range("A" + $row + ":D" + $row)

    # Copy format (style) down from the previous row, then overwrite values.
    $src.Copy($dst)

    $ws.Range("A" + $row).Value = $serial
    $ws.Range("B" + $row).Value = 0
    $ws.Range("C" + $row).Value = 0
    $ws.Range("D" + $row).Value = 0
}

Write-Output "Added rows $startRow to $($startRow + $numRows - 1)"
